$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Monday results update: sample rows 13 and 14 (the 10,000th-sample
# milestone upload) need their collection dates corrected, which means the
# two rows trade places entirely - sample id, dates, every lab/field
# column. Swap the full row contents (A:CW) between row 13 and row 14.
#
# Using Range.Copy (rather than reading/writing .Value2 cell-by-cell) so
# that text-formatted values which merely look numeric (sample ids, lat/lon,
# county FIPS codes with leading zeros, ...) keep their original text type
# and formatting instead of being reinterpreted as numbers.

$lastCol = $ws.UsedRange.Columns.Count
$lastColLetter = $ws.Cells.Item(1, $lastCol).Address($false, $false) -replace '\d+$', ''

$row13Range = $ws.Range("A13:" + $lastColLetter + "13")
$row14Range = $ws.Range("A14:" + $lastColLetter + "14")

# Scratch holding area far below the used range of the sheet.
$scratchRow = 100
$scratchRange = $ws.Range("A" + $scratchRow + ":" + $lastColLetter + $scratchRow)

# 1) stash row 13 in the scratch area
$row13Range.Copy($scratchRange)

# 2) move row 14's contents into row 13
$row14Range.Copy($row13Range)

# 3) move the stashed original row 13 contents into row 14
$scratchRange.Copy($row14Range)

# 4) tidy up the scratch area
$scratchRange.Clear()
